$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(151979169, 1),
    @(726896511, 1),
    @(550338825, 1),
    @(764228034, 1),
    @(758146980, 1),
    @(171935222, 2),
    @(20645182, 3),
    @(79418400, 1),
    @(641239853, 1),
    @(110151679, 1),
    @(749067789, 1),
    @(189153073, 1),
    @(254395074, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
